# Regenerate save_data to use K instead of Strike#: recalculated K (column G)
# values for each row based on the updated std/mean computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(6, 6, 13, 4, 11, 7, 7, 6, 12, 6, 8, 5, 5, 7, 5, 6, 4, 2, 11, 5, 3, 4, 6, 4, 4, 9, 7, 5, 5, 7, 2, 2, 2)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $newK[$i]
}
